# Update predictor labels to reflect log-transformed variables.
# Mapping of old text -> new text, as described by the commit:
#   fixing marginal effects (now conditional effects) to show both parts
#   of hurdle, and contour graph for interaction

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "GDP (dollars per capita)"                 = "ln(GDP [dollars per capita])"
    "Tourism - Inbound (per capita)"           = "ln(Tourism - Inbound [per capita])"
    "ProMed Mentions (per capita)"              = "ln(ProMed Mentions [per capita])"
    "Migrant Population (per capita)"           = "ln(Migrant Population [per capita])"
    "AB Exports (dollars per capita)"           = "ln(AB Exports [dollars per capita])"
    "Publication Bias Index (per capita)"       = "ln(Publication Bias Index [per capita])"
    "Livestock AB Consumption (kg per capita)"  = "Livestock AB Consumption [kg per capita)"
    "Population"                                = "ln(Population)"
}

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($null -ne $val -and $replacements.ContainsKey([string]$val)) {
            $cell.Value = $replacements[[string]$val]
        }
    }
}
